$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Cells.Item(6, 2).Value = 0.55149442314395203
$ws.Cells.Item(6, 3).Value = 0.53829295921924003
$ws.Cells.Item(6, 6).Value = 0.56744945974206995
$ws.Cells.Item(7, 2).Value = 0.552823283373997
$ws.Cells.Item(7, 3).Value = 0.546893516904844
$ws.Cells.Item(7, 6).Value = 0.55547664691530096
$ws.Cells.Item(8, 2).Value = 0.55221767166260005
$ws.Cells.Item(8, 3).Value = 0.51186824677588005
$ws.Cells.Item(8, 6).Value = 0.55212181944928496
$ws.Cells.Item(9, 2).Value = 0.55813872429417899
$ws.Cells.Item(9, 3).Value = 0.53291216451725298
$ws.Cells.Item(9, 6).Value = 0.55419135587312596
$ws.Cells.Item(10, 2).Value = 0.55020041826420296
$ws.Cells.Item(10, 3).Value = 0.54557337051237298
$ws.Cells.Item(10, 6).Value = 0.55945887068665001
$ws.Cells.Item(11, 2).Value = 0.56074851864761199
$ws.Cells.Item(11, 3).Value = 0.52310038340885301
$ws.Cells.Item(11, 6).Value = 0.55544179156500495
$ws.Cells.Item(12, 2).Value = 0.53302980132450295
$ws.Cells.Item(12, 3).Value = 0.52971418612757004
$ws.Cells.Item(12, 6).Value = 0.56010805158591803
$ws.Cells.Item(13, 2).Value = 0.55282764029278497
$ws.Cells.Item(13, 3).Value = 0.520473161380271
$ws.Cells.Item(13, 6).Value = 0.56471767166260001
$ws.Cells.Item(14, 2).Value = 0.53112147089578199
$ws.Cells.Item(14, 3).Value = 0.52971854304635702
$ws.Cells.Item(14, 6).Value = 0.55087138375740596
$ws.Cells.Item(15, 2).Value = 0.54819623562216802
$ws.Cells.Item(15, 3).Value = 0.55152927849424804
$ws.Cells.Item(15, 6).Value = 0.55810386894388198
$ws.Cells.Item(16, 2).Value = 0.55748954339491097
$ws.Cells.Item(16, 3).Value = 0.52514377831997205
$ws.Cells.Item(16, 6).Value = 0.55813872429417899
$ws.Cells.Item(17, 2).Value = 0.54420529801324502
$ws.Cells.Item(17, 3).Value = 0.54020564656674797
$ws.Cells.Item(17, 6).Value = 0.56141948414081499
$ws.Cells.Item(18, 2).Value = 0.55089752527012803
$ws.Cells.Item(18, 3).Value = 0.552897350993377
$ws.Cells.Item(18, 6).Value = 0.55555942837225503
$ws.Cells.Item(19, 2).Value = 0.55223509933774795
$ws.Cells.Item(19, 3).Value = 0.53032851167654205
$ws.Cells.Item(19, 6).Value = 0.56012547926106604
$ws.Cells.Item(20, 2).Value = 0.55221331474381297
$ws.Cells.Item(20, 3).Value = 0.54758191007319601
$ws.Cells.Item(20, 6).Value = 0.56209916347159194
$ws.Cells.Item(21, 2).Value = 0.53758278145695304
$ws.Cells.Item(21, 3).Value = 0.53567009410944499
$ws.Cells.Item(21, 6).Value = 0.56478738236319204
$ws.Cells.Item(22, 2).Value = 0.54358225862669896
$ws.Cells.Item(22, 3).Value = 0.55353346113628399
$ws.Cells.Item(22, 6).Value = 0.56671749738584798
$ws.Cells.Item(23, 2).Value = 0.55478389682816298
$ws.Cells.Item(23, 3).Value = 0.53499477169745502
$ws.Cells.Item(23, 6).Value = 0.55417828511676503
$ws.Cells.Item(24, 2).Value = 0.54823980481003798
$ws.Cells.Item(24, 3).Value = 0.55210003485534997
$ws.Cells.Item(24, 6).Value = 0.55939351690484496
$ws.Cells.Item(25, 2).Value = 0.54237539212269004
$ws.Cells.Item(25, 3).Value = 0.53037643778319898
$ws.Cells.Item(25, 6).Value = 0.54887591495294497
$ws.Cells.Item(26, 2).Value = 0.54297228999651403
$ws.Cells.Item(26, 3).Value = 0.54628790519344705
$ws.Cells.Item(26, 6).Value = 0.56145869640989898
$ws.Cells.Item(27, 2).Value = 0.55687957476472605
$ws.Cells.Item(27, 3).Value = 0.52171924015336302
$ws.Cells.Item(27, 6).Value = 0.55152056465667398
$ws.Cells.Item(28, 2).Value = 0.55940658766120599
$ws.Cells.Item(28, 3).Value = 0.54556901359358601
$ws.Cells.Item(28, 6).Value = 0.55883147438131697
$ws.Cells.Item(29, 2).Value = 0.53829295921924003
$ws.Cells.Item(29, 3).Value = 0.52711310561171099
$ws.Cells.Item(29, 6).Value = 0.55084524224468401
$ws.Cells.Item(30, 2).Value = 0.54495904496340197
$ws.Cells.Item(30, 3).Value = 0.53368333914255806
$ws.Cells.Item(30, 6).Value = 0.56006448239804796
$ws.Cells.Item(31, 2).Value = 0.56079644475426904
$ws.Cells.Item(31, 3).Value = 0.522486057859881
$ws.Cells.Item(31, 6).Value = 0.54959044963401804
$ws.Cells.Item(32, 2).Value = 0.550156849076333
$ws.Cells.Item(32, 3).Value = 0.56013855001742696
$ws.Cells.Item(32, 6).Value = 0.55877919135587295
$ws.Cells.Item(33, 2).Value = 0.53699459742070399
$ws.Cells.Item(33, 3).Value = 0.53231090972464201
$ws.Cells.Item(33, 6).Value = 0.55614761240850397
$ws.Cells.Item(34, 2).Value = 0.5475470547229
$ws.Cells.Item(34, 3).Value = 0.51251742767514796
$ws.Cells.Item(34, 6).Value = 0.54946409898919402
$ws.Cells.Item(35, 2).Value = 0.53757842453816596
$ws.Cells.Item(35, 3).Value = 0.55025270128964798
$ws.Cells.Item(35, 6).Value = 0.55614325548971699

# Scroll/selection change to match target view
$ws.Application.Goto($ws.Range("A21"), $false)
$ws.Range("H36").Select()
